# Module - Transactions_Jewel (Cash & Transfer) - Completed
#
# Adds a new worksheet "Transactions_Jewel" (modelled after the existing
# "GeneralOpening_SuspenseAsset" sheet) at the end of the workbook, and
# clears the tab-selection / cell-selection state that used to be on the
# previously-last sheet so the new sheet becomes the active one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Locate the current last sheet (GeneralOpening_SuspenseAsset) - this
#    is both the template for formatting and the sheet the new tab gets
#    inserted after.
# ---------------------------------------------------------------------
$sheetCount  = $wb.Worksheets.Count
$templateWs  = $wb.Worksheets.Item($sheetCount)

# ---------------------------------------------------------------------
# 2. Insert the new worksheet right after the template sheet and name it.
# ---------------------------------------------------------------------
$newWs = $wb.Worksheets.Add($null, $templateWs)
$newWs.Name = "Transactions_Jewel"

# ---------------------------------------------------------------------
# 3. Header row (row 1) - copy formatting from the template header cells,
#    then write the header text.
# ---------------------------------------------------------------------
$headerCols = @("A", "B", "C", "D", "E", "F", "G")
foreach ($col in $headerCols) {
    $templateWs.Range($col + "1").Copy()
    $newWs.Range($col + "1").PasteSpecial(-4122)  # xlPasteFormats
}

$newWs.Range("A1").Value = "TestScenario"
$newWs.Range("B1").Value = "Run"
$newWs.Range("C1").Value = "pcRegFormName"
$newWs.Range("D1").Value = "pcRegFormPcName"
$newWs.Range("E1").Value = "amount"
$newWs.Range("F1").Value = "remark"
$newWs.Range("G1").Value = "accNum"

# ---------------------------------------------------------------------
# 4. Data row (row 2) - copy formatting cell-by-cell (styles differ across
#    the row; E2 and G2 stay on the sheet's default/general style).
# ---------------------------------------------------------------------
$dataCols = @("A", "B", "C", "D", "F")
foreach ($col in $dataCols) {
    $templateWs.Range($col + "2").Copy()
    $newWs.Range($col + "2").PasteSpecial(-4122)  # xlPasteFormats
}

$newWs.Range("A2").Value = "Transactions_Jewel"
$newWs.Range("B2").Value = "Yes"
$newWs.Range("C2").Value = "qwerty"
$newWs.Range("D2").Value = "zxcvb"
$newWs.Range("E2").Value = 2000
$newWs.Range("F2").Value = "abcd"
$newWs.Range("G2").Value = 3

# ---------------------------------------------------------------------
# 5. Row heights for the new sheet (both rows are taller than default).
# ---------------------------------------------------------------------
$newWs.Rows.Item(1).RowHeight = 45
$newWs.Rows.Item(2).RowHeight = 45

# ---------------------------------------------------------------------
# 6. Selection / active-tab bookkeeping.
#    - the old last sheet loses its selection (J9) in favour of a block
#      selection over its header+data rows;
#    - the new sheet becomes selected / active with its own cursor cell.
# ---------------------------------------------------------------------
$templateWs.Range("A1:D2").Select()

$newWs.Activate()
$newWs.Range("L7").Select()

$wb.Save()
